$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.137.05"
$ws.Range("E2").Value = "  +3.01%  "
$ws.Range("D3").Value = "1.735.95"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("D4").Formula = "'0.9980"
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Formula = "'240.33"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").Formula = "'0.9975"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Formula = "'0.4801"
$ws.Range("E7").Value = "  -1.75%  "
$ws.Range("D8").Formula = "'0.2591"
$ws.Range("E8").Value = "  +0.55%  "
$ws.Range("D9").Formula = "'0.06142"
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("D10").Value = "1.720.83"
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("D11").Formula = "'16.03"
$ws.Range("E11").Value = "  +3.38%  "
$ws.Range("D12").Formula = "'0.06928"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Formula = "'0.6018"
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("D14").Formula = "'4.434"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").Formula = "'76.78"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").Formula = "'0.9995"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").Value = "27.110.02"
$ws.Range("E17").Value = "  +3.31%  "
$ws.Range("D18").Formula = "'0.9974"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("D19").Formula = "'0.000007050"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").Formula = "'11.38"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("D21").Value = "1.941.72"
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").Formula = "'4.401"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").Formula = "'8.379"
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").Formula = "'5.101"
$ws.Range("E24").Value = "  +1.50%  "
$ws.Range("D25").Formula = "'141.72"
$ws.Range("E25").Value = "  +3.99%  "
$ws.Range("D26").Formula = "'15.25"
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("D27").Formula = "'1.814"
$ws.Range("E27").Value = "  +5.23%  "
$ws.Range("D28").Formula = "'106.83"
$ws.Range("E28").Value = "  +1.32%  "
$ws.Range("D29").Formula = "'1.375"
$ws.Range("E29").Value = "  -1.10%  "
$ws.Range("D30").Formula = "'3.940"
$ws.Range("E30").Value = "  +1.53%  "
$ws.Range("D31").Formula = "'0.07924"
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").Formula = "'3.664"
$ws.Range("E32").Value = "  +1.94%  "
$ws.Range("D33").Formula = "'0.04723"
$ws.Range("E33").Value = "  +6.26%  "
$ws.Range("D34").Formula = "'2.593"
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("D35").Formula = "'1.010"
$ws.Range("E35").Value = "  +2.05%  "
$ws.Range("D36").Formula = "'0.6164"
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("D37").Formula = "'0.9216"
$ws.Range("E37").Value = "  -2.47%  "
$ws.Range("D38").Formula = "'2.536"
$ws.Range("E38").Value = "  +7.12%  "
$ws.Range("D39").Formula = "'2.014"
$ws.Range("E39").Value = "  +0.67%  "
$ws.Range("D40").Formula = "'0.9975"
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("D41").Formula = "'5.689"
$ws.Range("E41").Value = "  +6.02%  "
$ws.Range("D42").Formula = "'0.01487"
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("D43").Formula = "'99.03"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("D44").Formula = "'0.3823"
$ws.Range("E44").Value = "  +0.67%  "
$ws.Range("D45").Formula = "'6.835"
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("D46").Formula = "'0.1151"
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("D47").Formula = "'0.05345"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Formula = "'7.782"
$ws.Range("E48").Value = "  +1.11%  "
$ws.Range("D49").Formula = "'29.89"
$ws.Range("E49").Value = "  -1.62%  "
$ws.Range("D50").Formula = "'1.240"
$ws.Range("E50").Value = "  +3.32%  "
$ws.Range("D51").Formula = "'50.99"
$ws.Range("E51").Value = "  -0.22%  "
